$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S1")

# --- Window / view bookkeeping -------------------------------------------------
$wb.Windows.Item(1).Height = 14060

# Pane / selection on the S1 sheet view
$ws.Activate()
$ws.Range("H2").Select()
$ActiveWindow = $excel.ActiveWindow
$ActiveWindow.FreezePanes = $false
$ActiveWindow.SplitColumn = 3
$ActiveWindow.SplitRow = 1
$ActiveWindow.FreezePanes = $true
$ws.Range("O31").Select()

# --- New cells (value 1), copying formatting from an existing cell with the ---
# --- same visual style so borders / shading match --------------------------
$ws.Range("N24").Copy($ws.Range("N2"))
$ws.Range("N24").Copy($ws.Range("N3"))
$ws.Range("N24").Copy($ws.Range("O4"))
$ws.Range("N18").Copy($ws.Range("O5"))
$ws.Range("N28").Copy($ws.Range("O6"))
$ws.Range("N28").Copy($ws.Range("N7"))
$ws.Range("N28").Copy($ws.Range("O7"))
$ws.Range("N28").Copy($ws.Range("N8"))
$ws.Range("N28").Copy($ws.Range("O8"))
$ws.Range("N18").Copy($ws.Range("O9"))
$ws.Range("N18").Copy($ws.Range("O11"))
$ws.Range("N18").Copy($ws.Range("O12"))
$ws.Range("N18").Copy($ws.Range("O13"))
$ws.Range("N14").Copy($ws.Range("O14"))
$ws.Range("N15").Copy($ws.Range("O15"))
$ws.Range("N16").Copy($ws.Range("O16"))
$ws.Range("N16").Copy($ws.Range("O23"))
$ws.Range("N24").Copy($ws.Range("N25"))
$ws.Range("N28").Copy($ws.Range("O30"))

$ws.Range("A1").Select()
